# Appends sprint S09/S10 auth task rows (A1:I77 -> A1:I93) to Sheet1,
# matching the "S09/G01: Authentication backend (users, passwords, sessions)" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 78; A = "S09"; B = "G01"; C = "Authentication backend (users, passwords, sessions)"; D = "S09_G01_TB001"; E = "Design auth model (users table, roles, password hashing, session strategy) and update PRD/ERD accordingly."; F = "Chose a simple PBKDF2-SHA256 password hashing scheme and HMAC-signed session tokens stored in an HTTP-only cookie; no external auth libraries were introduced."; G = "implemented"; H = "Auth helpers live in app.core.auth and are shared between the runtime and Alembic migration used to seed the default admin."; I = "Document recommended auth-related environment variables (ST_CRYPTO_KEY for signing) and consider rotating secrets in a later sprint." },
    @{ Row = 79; A = "S09"; B = "G01"; C = "Authentication backend (users, passwords, sessions)"; D = "S09_G01_TB002"; E = "Implement SQLAlchemy User model + Alembic migration, seeding a default admin user (admin/admin) with a secure password hash."; F = "User ORM model and Alembic migration 0005_add_users.py added a users table with username/password_hash/role and timestamps, and seeded a default ADMIN user in the migration only."; G = "implemented"; H = "Tests create users via the API and models rather than relying on the seeded admin, since pytest uses metadata create_all instead of running migrations."; I = "When running Alembic migrations against an existing DB, verify that the seeded admin user is created only once and can be disabled or renamed later from the UI." },
    @{ Row = 80; A = "S09"; B = "G01"; C = "Authentication backend (users, passwords, sessions)"; D = "S09_G01_TB003"; E = "Implement backend auth APIs (register, login, logout, change password, current user) with JWT/cookie sessions and pytest coverage."; F = "Implemented minimal auth endpoints (register, login, logout, change-password, me) under /api/auth with cookie-based sessions; no route-level authorization changes were made yet."; G = "implemented"; H = "S09/G03 will attach these auth primitives to existing admin routers and refine role-based authorization; current APIs remain backwards compatible."; I = "Extend tests in later sprints to cover edge cases such as expired tokens and corrupted cookies once frontend auth flows are in place." },
    @{ Row = 81; A = "S09"; B = "G02"; C = "Frontend auth flows and landing layout"; D = "S09_G02_TF001"; E = "Create public login/register routes with a right-aligned auth panel and a marketing/hero area on the left describing SigmaTrader benefits."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 82; A = "S09"; B = "G02"; C = "Frontend auth flows and landing layout"; D = "S09_G02_TF002"; E = "Implement signup and login forms wired to auth APIs, including validation and error messaging."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 83; A = "S09"; B = "G02"; C = "Frontend auth flows and landing layout"; D = "S09_G02_TF003"; E = "Protect app routes behind login and add a top-right user menu (username, profile, change password, logout)."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 84; A = "S09"; B = "G03"; C = "Authorization and integration with existing admin features"; D = "S09_G03_TB001"; E = "Replace or augment HTTP Basic admin protection with role-based user auth, mapping ADMIN role to existing admin-only APIs."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 85; A = "S09"; B = "G03"; C = "Authorization and integration with existing admin features"; D = "S09_G03_TB002"; E = "Ensure broker config, risk settings, analytics, and system events remain restricted to admin users while normal traders can view their own data."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 86; A = "S09"; B = "G03"; C = "Authorization and integration with existing admin features"; D = "S09_G03_TB003"; E = "Define dev-mode behaviour (optional auth bypass) and verify TradingView webhook and Zerodha connect flows remain compatible."; G = "pending"; I = "Planned for S09; depends on existing Zerodha integration and risk/analytics APIs being stable." },
    @{ Row = 87; A = "S10"; B = "G01"; C = "Auth security refinements (rate limiting, password reset)"; D = "S10_G01_TB001"; E = "Add basic rate limiting / lockout behaviour on login to reduce brute-force attempts (e.g., small delay or temporary block after repeated failures)."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 88; A = "S10"; B = "G01"; C = "Auth security refinements (rate limiting, password reset)"; D = "S10_G01_TB002"; E = "Implement password reset flows: change-password for users plus an admin-only endpoint to reset another user's password."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 89; A = "S10"; B = "G02"; C = "Auth observability and audit logging"; D = "S10_G02_TB001"; E = "Record login, logout, and password-change events into system_events and surface them in the System Events UI with appropriate filters."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 90; A = "S10"; B = "G02"; C = "Auth observability and audit logging"; D = "S10_G02_TB002"; E = "Add optional notifications or banners when suspicious auth activity is detected (e.g., repeated failures)."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 91; A = "S10"; B = "G03"; C = "Roles and user experience enhancements"; D = "S10_G03_TB001"; E = "Introduce additional roles (e.g., VIEW_ONLY) and adjust API/UI permissions so that view-only users cannot modify risk, strategies, or broker config."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 92; A = "S10"; B = "G03"; C = "Roles and user experience enhancements"; D = "S10_G03_TB002"; E = "Add per-user preferences (e.g., default landing page, theme choice) stored in the DB and applied in the frontend."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." },
    @{ Row = 93; A = "S10"; B = "G04"; C = "Future multi-broker/multi-account design (auth-aware)"; D = "S10_G04_TB001"; E = "Design how users map to brokers/accounts (e.g., single Zerodha account vs per-user broker credentials) and document the migration path from the current single-account model."; G = "pending"; I = "Planned refinements for post-S09 auth, focused on security, auditing, and UX." }
)

$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8; I = 9 }

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
        if ($rowData.ContainsKey($col)) {
            $cell = $ws.Cells.Item($r, $colIndex[$col])
            $cell.Value = $rowData[$col]
            $cell.Style = "Normal"
            $cell.VerticalAlignment = -4107
        }
    }
}

